$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 896361
$ws.Cells.Item(3, 3).Value = 3570046
$ws.Cells.Item(4, 3).Value = 8031321
$ws.Cells.Item(5, 3).Value = 14443122
$ws.Cells.Item(6, 3).Value = 21608873
$ws.Cells.Item(7, 3).Value = 31548381
$ws.Cells.Item(8, 3).Value = 43043691
$ws.Cells.Item(9, 3).Value = 55778767
$ws.Cells.Item(10, 3).Value = 68954324
$ws.Cells.Item(11, 3).Value = 87387447
$ws.Cells.Item(12, 3).Value = 105738032
$ws.Cells.Item(13, 3).Value = 122630750
$ws.Cells.Item(14, 3).Value = 149761477
$ws.Cells.Item(15, 3).Value = 170666861
$ws.Cells.Item(16, 3).Value = 196267445
$ws.Cells.Item(17, 3).Value = 222976935
$ws.Cells.Item(18, 3).Value = 252088550
$ws.Cells.Item(19, 3).Value = 282288114
$ws.Cells.Item(20, 3).Value = 308249638
$ws.Cells.Item(21, 3).Value = 347867402
$ws.Cells.Item(22, 3).Value = 382815197
$ws.Cells.Item(23, 3).Value = 413627587
$ws.Cells.Item(24, 3).Value = 460594138
$ws.Cells.Item(25, 3).Value = 504784367
$ws.Cells.Item(26, 3).Value = 546119157
$ws.Cells.Item(27, 3).Value = 588423436
$ws.Cells.Item(28, 3).Value = 626265157
$ws.Cells.Item(29, 3).Value = 679315499
$ws.Cells.Item(30, 3).Value = 732459496
$ws.Cells.Item(31, 3).Value = 767312779
$ws.Cells.Item(32, 3).Value = 837357190
$ws.Cells.Item(33, 3).Value = 893678649
$ws.Cells.Item(34, 3).Value = 938133167
$ws.Cells.Item(35, 3).Value = 1004456571
$ws.Cells.Item(36, 3).Value = 1063519284
$ws.Cells.Item(37, 3).Value = 1125172719
$ws.Cells.Item(38, 3).Value = 1188437121
$ws.Cells.Item(39, 3).Value = 1256221815
$ws.Cells.Item(40, 3).Value = 1320900038
$ws.Cells.Item(41, 3).Value = 1398677697
$ws.Cells.Item(42, 3).Value = 1445158007
$ws.Cells.Item(43, 3).Value = 1532050232
$ws.Cells.Item(44, 3).Value = 1603760351
$ws.Cells.Item(45, 3).Value = 1651725045
$ws.Cells.Item(46, 3).Value = 1729779395
$ws.Cells.Item(47, 3).Value = 1807971450
$ws.Cells.Item(48, 3).Value = 1887325010
$ws.Cells.Item(49, 3).Value = 1962021429
$ws.Cells.Item(50, 3).Value = 2048261056
$ws.Cells.Item(51, 3).Value = 2130755299
$ws.Cells.Item(52, 3).Value = 2230060577
$ws.Cells.Item(53, 3).Value = 2310836938
$ws.Cells.Item(54, 3).Value = 2388531205
$ws.Cells.Item(55, 3).Value = 2490529841
$ws.Cells.Item(56, 3).Value = 2592373239
$ws.Cells.Item(57, 3).Value = 2675116375
$ws.Cells.Item(58, 3).Value = 2764769955
$ws.Cells.Item(59, 3).Value = 2866622333
$ws.Cells.Item(60, 3).Value = 2967947843
$ws.Cells.Item(61, 3).Value = 3060439755
$ws.Cells.Item(62, 3).Value = 3171168497
$ws.Cells.Item(63, 3).Value = 3273232992
$ws.Cells.Item(64, 3).Value = 3371293372
$ws.Cells.Item(65, 3).Value = 3503242678
$ws.Cells.Item(66, 3).Value = 3609342078
$ws.Cells.Item(67, 3).Value = 3711430948
$ws.Cells.Item(68, 3).Value = 3834827197
$ws.Cells.Item(69, 3).Value = 3944926433
$ws.Cells.Item(70, 3).Value = 4053865448
$ws.Cells.Item(71, 3).Value = 4174434910
$ws.Cells.Item(72, 3).Value = 4290173178
$ws.Cells.Item(73, 3).Value = 4414925941
$ws.Cells.Item(74, 3).Value = 4542305348
$ws.Cells.Item(75, 3).Value = 4664282217
$ws.Cells.Item(76, 3).Value = 4788798916
$ws.Cells.Item(77, 3).Value = 4916803123
$ws.Cells.Item(78, 3).Value = 5048039960
$ws.Cells.Item(79, 3).Value = 5197523034
$ws.Cells.Item(80, 3).Value = 5312369646
$ws.Cells.Item(81, 3).Value = 5453159820
$ws.Cells.Item(82, 3).Value = 5581098170
$ws.Cells.Item(83, 3).Value = 5732867904
$ws.Cells.Item(84, 3).Value = 5857242195
$ws.Cells.Item(85, 3).Value = 6153465806
$ws.Cells.Item(86, 3).Value = 6147170756
$ws.Cells.Item(87, 3).Value = 6287211682
$ws.Cells.Item(88, 3).Value = 6444198782
$ws.Cells.Item(89, 3).Value = 6600765926
$ws.Cells.Item(90, 3).Value = 6745422939
$ws.Cells.Item(91, 3).Value = 6909314922
$ws.Cells.Item(92, 3).Value = 7049804456
$ws.Cells.Item(93, 3).Value = 7207017784
$ws.Cells.Item(94, 3).Value = 7367514203
$ws.Cells.Item(95, 3).Value = 7512111343
$ws.Cells.Item(96, 3).Value = 7682642211
$ws.Cells.Item(97, 3).Value = 7845687869
$ws.Cells.Item(98, 3).Value = 8036432606
$ws.Cells.Item(99, 3).Value = 8183771000
$ws.Cells.Item(100, 3).Value = 8339085122
$ws.Cells.Item(101, 3).Value = 8515346538
